$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F9").Value = 4
$ws.Range("F11").Value = -2
$ws.Range("F34").Value = -1
$ws.Range("F35").Value = 4
$ws.Range("F39").Value = -4
$ws.Range("F40").Value = -3
$ws.Range("F41").Value = 4
$ws.Range("F42").Value = 4
$ws.Range("F43").Value = -1
$ws.Range("F44").Value = -1
$ws.Range("F45").Value = 4
